$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values (registration entry changed)
$ws.Range("B3").Value = "Luaniinha23"
$ws.Range("C3").Value = "Batman_Bolado"
$ws.Range("D3").Value = "Batman Bolado"

# Add new row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ieko"
$ws.Range("C4").Value = "leko"
$ws.Range("D4").Value = "leko"
